$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("NetSuite") is rewritten in place to describe the Salesforce system instead.
$ws.Range("A2").Value = "Salesforce"
$ws.Range("B2").Value = "This is a cloud-based CRM platform providing sales, marketing, and service functionalities. It supports lead management, sales tracking, and integration with financial systems for data transfer."
$ws.Range("C2").Value = "The client utilizes this system as the primary CRM solution, tracking clients, managing sales opportunities, and transferring sales data to the ERP system."
$ws.Range("D2").Value = "The Salesforce Administration Team, led by Vincent Washburn (Senior Salesforce Engineer)."
$ws.Range("E2").Value = "Access provisioning is triggered by HR updates in Workday for new hires, which automatically create a ticket in ServiceNow. The ticket includes necessary role information and requires manager confirmation before access is provisioned by the IT team. Additional access requests are handled informally through direct communication with the IT team, bypassing a formal ticketing process. Role changes were not explicitly addressed."
$ws.Range("F2").Value = "When HR enters a termination date in Workday, an automated process creates a ServiceNow ticket, notifying system owners to remove access within three business days and requiring confirmation from the system owner before the ticket is closed."
$ws.Range("G2").Value = "Access is configured using a role-based model, where roles and permissions are determined by job titles in Workday and provisioned automatically through ServiceNow, with additional access manually granted by administrators upon informal requests."
$ws.Range("H2").Value = "Yes; system administrators can modify roles, configurations, and integrations, and they also have access to the back-end code."
$ws.Range("I2").Value = "N/A - This information was not discussed in the walkthrough meeting transcript."
$ws.Range("J2").Value = "Yes; Users with privileged access include Vincent Washburn and Jeremy, who have administrative access to the Salesforce application."
$ws.Range("K2").Value = "No; there are no interactive system accounts within the system."
$ws.Range("L2").Value = "N/A - No Interactive System Accounts"
$ws.Range("N2").Value = "Yes; Management conducts quarterly access reviews. Compliance initiates the process, and user access is reviewed for appropriateness during each cycle."
$ws.Range("O2").Value = "The system maintains logs of user activities. Logs are retained within Salesforce and include administrative user actions."
$ws.Range("P2").Value = "No; the system logs are not reviewed unless there is an incident"
$ws.Range("Q2").Value = "The system uses a hybrid authentication approach where regular users authenticate through Okta SSO, while administrators maintain direct application credentials as a backup authentication path."
$ws.Range("R2").Value = "Yes; The client periodically reviews authentication configurations to ensure alignment with security standards."
$ws.Range("S2").Value = "Management can perform configuration changes to adjust system settings, update workflows, and modify custom code within Salesforce. These changes are managed by Vincent Washburn and his colleague Jeremy, who have the ability to develop and deploy updates in a sandbox environment before moving to production."
$ws.Range("T2").Value = "The following roles have change capabilities: Senior Salesforce Engineer (administrative and configuration changes) and Salesforce Engineer (administrative and configuration changes). These roles are assigned to Vincent Washburn and his colleague Jeremy, who are the only individuals with administrative access to Salesforce and are authorized to make changes to the application."
$ws.Range("U2").Value = "Management maintains a separate sandbox environment for implementing and testing changes."
$ws.Range("V2").Value = "The client's change management process includes the following steps: 1) Change request submission via JIRA ticket created by Vincent or Jeremy, 2) Development in the sandbox environment by Vincent or Jeremy, 3) Testing in the sandbox environment coordinated by Jeremy with relevant stakeholders, 4) Approval captured within the JIRA ticket by the Change Advisory Board (CAB) before deployment, 5) Deployment to production performed by Vincent or Jeremy without strict segregation of duties between development and deployment."
$ws.Range("W2").Value = "Yes, the vendor manages updates, patches, and bug fixes for the SaaS system, but the specific process for deploying these changes to production was not discussed."
$ws.Range("Y2").Value = "No periodic review of changes is performed."
$ws.Range("Z2").Value = "Yes, the system has automated jobs and interfaces with NetSuite. These jobs transfer sales information for closed sales, including numbers and client data, using a tool that pulls data daily and transforms it for compatibility."
$ws.Range("AA2").Value = "Jobs are managed using Salesforce's native functionality, NetSuite's native functionality, and Workato. These tools provide scheduling and monitoring capabilities for automated jobs and interfaces."
$ws.Range("AB2").Value = "Job failures are handled through a defined process: detection is performed manually by the accounting team, notification occurs when the accounting team informs the IT team of issues, resolution involves reviewing Workato logs and rerunning the job or opening a ticket for complex problems."
$ws.Range("AE2").Value = "N/A - All backup types and strategies are managed by the vendor as part of their SaaS service."
$ws.Range("AF2").Value = "Backup failures are addressed through monitoring by the accounting team, notification to the IT team, review of Workato logs to identify errors, rerunning the job to resolve issues, and opening a ticket for complex problems with resolution logged in the ticket."
$ws.Range("AG2").Value = "Yes, management performs annual SOC 1 Type 2 report reviews. The review process includes evaluating complementary user entity controls and subservice providers and is documented in the compliance team's review template."

# Row 3 ("Microsoft Azure") is removed entirely; remaining rows shift up so the
# used range collapses from A1:AG3 down to A1:AG2.
$ws.Rows(3).Delete()
